$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text (preventing Excel from
# auto-converting numeric-looking strings like "1.00" or "0.999"
# into real numbers), while keeping the cell on the default/Normal
# style so no stray number-format style is left behind.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "63.330.14"
Set-TextValue $ws.Range("E2") "  -3.83%  "
Set-TextValue $ws.Range("D3") "3.289.54"
Set-TextValue $ws.Range("E3") "  -6.26%  "
Set-TextValue $ws.Range("D4") "0.999"
Set-TextValue $ws.Range("E4") "  -0.18%  "
Set-TextValue $ws.Range("D5") "547.74"
Set-TextValue $ws.Range("E5") "  -4.53%  "
Set-TextValue $ws.Range("D6") "171.71"
Set-TextValue $ws.Range("E6") "  -3.82%  "
Set-TextValue $ws.Range("D7") "0.607"
Set-TextValue $ws.Range("E7") "  -4.22%  "
Set-TextValue $ws.Range("E8") "  +0.04%  "
Set-TextValue $ws.Range("D9") "3.286.74"
Set-TextValue $ws.Range("E9") "  -6.17%  "
Set-TextValue $ws.Range("D10") "0.615"
Set-TextValue $ws.Range("E10") "  -2.88%  "
Set-TextValue $ws.Range("D11") "0.158"
Set-TextValue $ws.Range("E11") "  -0.35%  "
Set-TextValue $ws.Range("D12") "53.21"
Set-TextValue $ws.Range("E12") "  -2.70%  "
Set-TextValue $ws.Range("D13") "0.0000268"
Set-TextValue $ws.Range("E13") "  -1.85%  "
Set-TextValue $ws.Range("D14") "8.92"
Set-TextValue $ws.Range("E14") "  -3.29%  "
Set-TextValue $ws.Range("D15") "3.814.14"
Set-TextValue $ws.Range("E15") "  -6.31%  "
Set-TextValue $ws.Range("D16") "18.03"
Set-TextValue $ws.Range("E16") "  -1.47%  "
Set-TextValue $ws.Range("D17") "0.117"
Set-TextValue $ws.Range("E17") "  -3.82%  "
Set-TextValue $ws.Range("D18") "3.281.89"
Set-TextValue $ws.Range("E18") "  -6.74%  "
Set-TextValue $ws.Range("D19") "11.65"
Set-TextValue $ws.Range("E19") "  -3.78%  "
Set-TextValue $ws.Range("D20") "63.141.06"
Set-TextValue $ws.Range("E20") "  -4.28%  "
Set-TextValue $ws.Range("D21") "0.961"
Set-TextValue $ws.Range("E21") "  -4.44%  "
Set-TextValue $ws.Range("D22") "421.94"
Set-TextValue $ws.Range("E22") "  +1.68%  "
Set-TextValue $ws.Range("D23") "4.59"
Set-TextValue $ws.Range("E23") "  +4.66%  "
Set-TextValue $ws.Range("D24") "4.05"
Set-TextValue $ws.Range("E24") "  -3.10%  "
Set-TextValue $ws.Range("D25") "83.62"
Set-TextValue $ws.Range("E25") "  -2.14%  "
Set-TextValue $ws.Range("D26") "13.00"
Set-TextValue $ws.Range("E26") "  +1.78%  "
Set-TextValue $ws.Range("D27") "10.54"
Set-TextValue $ws.Range("E27") "  -3.50%  "
Set-TextValue $ws.Range("D28") "2.79"
Set-TextValue $ws.Range("E28") "  -2.11%  "
Set-TextValue $ws.Range("D29") "8.58"
Set-TextValue $ws.Range("E29") "  -4.59%  "
Set-TextValue $ws.Range("D30") "29.26"
Set-TextValue $ws.Range("E30") "  -3.50%  "
Set-TextValue $ws.Range("D31") "6.54"
Set-TextValue $ws.Range("E31") "  +2.11%  "
Set-TextValue $ws.Range("D32") "584.83"
Set-TextValue $ws.Range("E32") "  -5.79%  "
Set-TextValue $ws.Range("D33") "11.31"
Set-TextValue $ws.Range("E33") "  -2.93%  "
Set-TextValue $ws.Range("D34") "0.106"
Set-TextValue $ws.Range("E34") "  -4.04%  "
Set-TextValue $ws.Range("D35") "57.79"
Set-TextValue $ws.Range("E35") "  -3.14%  "
Set-TextValue $ws.Range("D36") "1.00"
Set-TextValue $ws.Range("E36") "  -0.02%  "
Set-TextValue $ws.Range("D37") "0.143"
Set-TextValue $ws.Range("E37") "  -6.66%  "
Set-TextValue $ws.Range("E38") "  +2.77%  "
Set-TextValue $ws.Range("D39") "35.02"
Set-TextValue $ws.Range("E39") "  -6.04%  "
Set-TextValue $ws.Range("D40") "0.0₃0738"
Set-TextValue $ws.Range("E40") "  -8.16%  "
Set-TextValue $ws.Range("D41") "0.361"
Set-TextValue $ws.Range("E41") "  -4.85%  "
Set-TextValue $ws.Range("D42") "0.998"
Set-TextValue $ws.Range("E42") "  -0.38%  "
Set-TextValue $ws.Range("D43") "3.064.43"
Set-TextValue $ws.Range("E43") "  -6.35%  "
Set-TextValue $ws.Range("E44") "  -4.50%  "
Set-TextValue $ws.Range("D45") "3.18"
Set-TextValue $ws.Range("E45") "  -2.70%  "
Set-TextValue $ws.Range("D46") "0.0400"
Set-TextValue $ws.Range("E46") "  -4.15%  "
Set-TextValue $ws.Range("D47") "2.41"
Set-TextValue $ws.Range("E47") "  -3.71%  "
Set-TextValue $ws.Range("E48") "  -3.05%  "
Set-TextValue $ws.Range("D49") "2.56"
Set-TextValue $ws.Range("E49") "  -5.61%  "
Set-TextValue $ws.Range("D50") "132.26"
Set-TextValue $ws.Range("E50") "  -4.56%  "
Set-TextValue $ws.Range("D51") "8.05"
Set-TextValue $ws.Range("E51") "  -4.75%  "
